# Apply cryptos list update (prices & volume %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '65.635.37'
$ws.Cells.Item(2, 5).Value = '  +1.61%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.476.74'
$ws.Cells.Item(3, 5).Value = '  +0.35%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '580.38'
$ws.Cells.Item(5, 5).Value = '  +0.21%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '161.55'
$ws.Cells.Item(6, 5).Value = '  +2.78%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.480.51'
$ws.Cells.Item(8, 5).Value = '  +0.29%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.584'
$ws.Cells.Item(9, 5).Value = '  +4.57%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.31'
$ws.Cells.Item(10, 5).Value = '  -3.85%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.126'
$ws.Cells.Item(11, 5).Value = '  +0.64%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.445'
$ws.Cells.Item(12, 5).Value = '  -0.56%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.082.67'
$ws.Cells.Item(13, 5).Value = '  +0.42%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.134'
$ws.Cells.Item(14, 5).Value = '  -1.57%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000196'
$ws.Cells.Item(15, 5).Value = '  -0.99%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '28.64'
$ws.Cells.Item(16, 5).Value = '  +3.07%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '65.610.25'
$ws.Cells.Item(17, 5).Value = '  +1.52%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.484.95'
$ws.Cells.Item(18, 5).Value = '  +0.45%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.46'
$ws.Cells.Item(19, 5).Value = '  +0.11%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.34'
$ws.Cells.Item(20, 5).Value = '  -0.27%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '392.68'
$ws.Cells.Item(21, 5).Value = '  -1.14%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '8.26'
$ws.Cells.Item(22, 5).Value = '  -3.73%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Polygon'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.551'
$ws.Cells.Item(23, 5).Value = '  +0.82%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '73.64'
$ws.Cells.Item(24, 5).Value = '  +0.89%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.03%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0000125'
$ws.Cells.Item(26, 5).Value = '  +2.91%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.63'
$ws.Cells.Item(27, 5).Value = '  +1.06%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.180'
$ws.Cells.Item(28, 5).Value = '  +0.08%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.29%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.43'
$ws.Cells.Item(30, 5).Value = '  +7.42%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.45'
$ws.Cells.Item(31, 5).Value = '  +4.51%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.06'
$ws.Cells.Item(32, 5).Value = '  +0.47%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.66'
$ws.Cells.Item(33, 5).Value = '  +0.20%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '23.75'
$ws.Cells.Item(34, 5).Value = '  -0.29%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '7.16'
$ws.Cells.Item(36, 5).Value = '  +2.66%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +2.30%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '161.46'
$ws.Cells.Item(38, 5).Value = '  +0.74%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.97'
$ws.Cells.Item(39, 5).Value = '  +4.53%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.051.87'
$ws.Cells.Item(40, 5).Value = '  +5.24%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0775'
$ws.Cells.Item(41, 5).Value = '  -1.25%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '27.12'
$ws.Cells.Item(42, 5).Value = '  -3.11%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.0323'
$ws.Cells.Item(43, 5).Value = '  +0.01%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.55'
$ws.Cells.Item(44, 5).Value = '  +2.75%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '42.77'
$ws.Cells.Item(45, 5).Value = '  +2.48%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.776'
$ws.Cells.Item(46, 5).Value = '  -0.65%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '25.90'
$ws.Cells.Item(47, 5).Value = '  +12.59%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.12'
$ws.Cells.Item(48, 5).Value = '  +2.56%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.26'
$ws.Cells.Item(49, 5).Value = '  +4.02%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '6.72'
$ws.Cells.Item(50, 5).Value = '  +2.49%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Bittensor'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '310.19'
$ws.Cells.Item(51, 5).Value = '  +4.25%  '
